$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two purchase-order numbers in column I
$ws.Range("I12").Value = 11520192875
$ws.Range("I13").Value = 11520174353

# Fill the header band (row 10, columns A:L) with sequential test values 0..11
$ws.Range("A10").Value = 0
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 6
$ws.Range("H10").Value = 7
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 9
$ws.Range("K10").Value = 10
$ws.Range("L10").Value = 11

# Add a new test entry below the table
$ws.Range("E14").Value = "dasdsa"

# Move the visible selection, matching the author's last position on screen
$ws.Range("F18").Select()
